$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "indice"
$ws.Range("B1").Value = "codigo"
$ws.Range("C1").Value = "descripcion"
$ws.Range("D1").Value = "marca"
$ws.Range("E1").Value = "cantidad"
$ws.Range("F1").Value = "precion sin igv"
$ws.Range("G1").Value = "precio"
$ws.Range("H1").Value = "Total_sin_igv"

# Copy header style (bold, centered, bordered) from A1 onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("D1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Range("A2").Value = 1

# B2 must stay text "3" (not get auto-converted to a number)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "mouse"
$ws.Range("D2").Value = " "
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 8.470000000000001
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 8.470000000000001

# --- Row 3 ---
$ws.Range("A3").Value = 2

# B3 must stay text "23" (not get auto-converted to a number)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "23"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").Value = "mouse 32"
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 8.470000000000001
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 16.95

# --- Remove old row 4 ---
$ws.Range("A4:E4").ClearContents()
